# Apply the "Saldo" workbook update:
#  - Row 2 (003301389 / EDMUNDO): Saldo 123698.6 -> 113000
#  - Delete row for 004334062 / MERG / 33855.83
#  - Delete row for 004693349 / CATARINE / 16044.94
#  - Delete row for 004805273 / CLISIA / 12345.94
#  - Replace row 004482102 / NATALIA / 6000 with 004335144 / EDMUNDO / 10698.6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update EDMUNDO's (003301389) balance on row 2.
$ws.Range("C2").Value = 113000

# 2) Delete the three obsolete rows entirely (bottom-up so row numbers
#    of not-yet-processed rows stay valid).
$ws.Range("A9").EntireRow.Delete()   # 004805273 / CLISIA / 12345.94
$ws.Range("A8").EntireRow.Delete()   # 004693349 / CATARINE / 16044.94
$ws.Range("A6").EntireRow.Delete()   # 004334062 / MERG / 33855.83

# 3) Replace the NATALIA row (now shifted up to row 9, since 3 rows
#    above it were deleted) with the new EDMUNDO entry. Force the
#    account number column to text so the leading zeros survive
#    (matches the other "Conta" cells, which are all text-typed).
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "004335144"
$ws.Range("B9").Value = "EDMUNDO"
$ws.Range("C9").Value = 10698.6

Write-Host "Saldo update applied: C2=$($ws.Range('C2').Value()), row9=$($ws.Range('A9').Value())/$($ws.Range('B9').Value())/$($ws.Range('C9').Value())"
